$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.574.08"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.862.44"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3897"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07898"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9735"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.811.00"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.729"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.948"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06912"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "28.606.85"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.339"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("D25").Value = "2.104.92"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.794"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.000"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09310"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9398"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.318"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.337"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.336"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05851"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.151"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5658"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.947"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1773"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07356"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5314"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.169"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.136"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.850"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.351"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("E51").Value = "  +0.15%  "
